$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so that numeric-looking
# values (e.g. "2.95", "65.74") are not auto-converted to numbers by Excel,
# matching the original inline-string cell content.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '37.345.09'
$ws.Range("E2").Value = '  +4.27%  '
$ws.Range("D3").Value = '2.043.23'
$ws.Range("E3").Value = '  +2.82%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '251.72'
$ws.Range("E5").Value = '  +2.79%  '
$ws.Range("D6").Value = '0.650'
$ws.Range("E6").Value = '  +1.97%  '
$ws.Range("D7").Value = '65.74'
$ws.Range("E7").Value = '  +11.13%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  +9.37%  '
$ws.Range("D10").Value = '59.41'
$ws.Range("E10").Value = '  +2.23%  '
$ws.Range("D11").Value = '0.0801'
$ws.Range("E11").Value = '  +8.59%  '
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").Value = '0.909'
$ws.Range("E13").Value = '  -3.48%  '
$ws.Range("D14").Value = '23.34'
$ws.Range("E14").Value = '  +22.57%  '
$ws.Range("D15").Value = '14.79'
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("D16").Value = '2.340.83'
$ws.Range("E16").Value = '  +2.75%  '
$ws.Range("D17").Value = '5.74'
$ws.Range("E17").Value = '  +8.07%  '
$ws.Range("D18").Value = '2.042.03'
$ws.Range("E18").Value = '  +2.62%  '
$ws.Range("D19").Value = '37.260.43'
$ws.Range("E19").Value = '  +4.29%  '
$ws.Range("D20").Value = '73.08'
$ws.Range("E20").Value = '  +2.19%  '
$ws.Range("D21").Value = '0.0₃0891'
$ws.Range("E21").Value = '  +5.34%  '
$ws.Range("D22").Value = '5.50'
$ws.Range("E22").Value = '  +6.07%  '
$ws.Range("D23").Value = '238.83'
$ws.Range("E23").Value = '  +2.68%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").Value = '2.36'
$ws.Range("E26").Value = '  +4.03%  '
$ws.Range("D27").Value = '10.00'
$ws.Range("E27").Value = '  +4.11%  '
$ws.Range("D28").Value = '161.34'
$ws.Range("E28").Value = '  -1.99%  '
$ws.Range("D29").Value = '20.02'
$ws.Range("E29").Value = '  +3.74%  '
$ws.Range("D30").Value = '0.130'
$ws.Range("E30").Value = '  +33.03%  '
$ws.Range("E31").Value = '  +2.91%  '
$ws.Range("D32").Value = '5.14'
$ws.Range("E32").Value = '  +4.63%  '
$ws.Range("D33").Value = '1.19'
$ws.Range("E33").Value = '  +5.44%  '
$ws.Range("D34").Value = '0.0626'
$ws.Range("E34").Value = '  +4.42%  '
$ws.Range("D35").Value = '4.65'
$ws.Range("E35").Value = '  +5.90%  '
$ws.Range("D36").Value = '6.40'
$ws.Range("E36").Value = '  +12.38%  '
$ws.Range("D37").Value = '2.37'
$ws.Range("E37").Value = '  -2.97%  '
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("E39").Value = '  +3.12%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '2.95'
$ws.Range("E40").Value = '  +27.81%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '1.29'
$ws.Range("E41").Value = '  +5.19%  '
$ws.Range("D42").Value = '0.102'
$ws.Range("E42").Value = '  +10.23%  '
$ws.Range("E43").Value = '  +5.68%  '
$ws.Range("D44").Value = '1.17'
$ws.Range("E44").Value = '  +6.03%  '
$ws.Range("D45").Value = '17.38'
$ws.Range("E45").Value = '  +5.96%  '
$ws.Range("D46").Value = '0.0219'
$ws.Range("E46").Value = '  +2.83%  '
$ws.Range("D47").Value = '95.43'
$ws.Range("E47").Value = '  +2.44%  '
$ws.Range("D48").Value = '7.83'
$ws.Range("E48").Value = '  +1.28%  '
$ws.Range("D49").Value = '1.395.92'
$ws.Range("E49").Value = '  +2.44%  '
$ws.Range("D50").Value = '2.92'
$ws.Range("E50").Value = '  +1.01%  '
$ws.Range("D51").Value = '47.09'
$ws.Range("E51").Value = '  +1.17%  '
